$p = $ppt.ActivePresentation

# Slide 11: "Linear regression is just a special case." textbox.
# Shrink the shape's height and fix the typo'd third paragraph.
$s11 = $p.Slides.Item(11)
$shp11 = $s11.Shapes.Item(3)
$shp11.Height = 2492990 / 12700
$para = $shp11.TextFrame.TextRange.Paragraphs(3)
$para.Runs(1).Text = "Just remove the hidden layer."

# Slide 7: add missing period at end of sentence.
$s7 = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(15)
$shp7.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Each node is a dot product of previous nodes in the layer before it."

# Slide 8: same fix, duplicated text box.
$s8 = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item(15)
$shp8.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Each node is a dot product of previous nodes in the layer before it."
